$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Solar thermal" water-heating row below the existing data
$ws.Range("A9").Value = "Solar thermal"
$ws.Range("B9").Value = "RES_CWH_SOLAR"
$ws.Range("C9").Value = "Water heating"

# Match the formatting (fill/border) used by the other data rows
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("C2").Select()
